# Workbook under edit: "Questionare on Regression testing.xlsx"
#
# 1. Rename the "LogError" sheet to "StatusSheet77618".
# 2. On "Sheet1", remove the stray "Status" header that was sitting in G1
#    (the sheet only uses columns A:F, rows 1:21 - clearing it also lets
#    the sheet's used-range/dimension shrink back down to A1:F21).

$wb = $excel.ActiveWorkbook

$logSheet = $wb.Worksheets.Item("LogError")
$logSheet.Name = "StatusSheet77618"

$mainSheet = $wb.Worksheets.Item("Sheet1")
$mainSheet.Range("G1").ClearContents()
